# Update Adv. Search Diagram
# - Merge the title's two runs into a single run of text.
# - Reposition/resize the "University Student (User)" actor rectangle.
# - Re-capitalise "student" -> "Student" (split into 3 runs).
# - Reposition/resize the dashed lifeline connector hanging off that rectangle.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title 3 (shape id 4): "Advanced Search " + "Scenario" -> "Advanced Search Scenario"
# (the two runs already concatenate to the same text, so re-assigning the
# identical string is a no-op for the engine; round-trip through a throwaway
# value first to force the two runs to collapse into a single run)
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "X"
$title.TextFrame.TextRange.Text = "Advanced Search Scenario"

# --- Rectangle 4 (shape id 5): move/resize, and fix actor label text
# (EMU -> point values picked so the round-trip back to EMU lands exactly
# on the target offsets from the diff: 1789041 / 3200402 EMU)
$rect = $s.Shapes.Item(2)
$rect.Left = 140.8694
$rect.Width = 252.0002

$tr = $rect.TextFrame.TextRange
$tr.Characters(1, 11).Text = "University "
$tr.Characters(12, 8).Text = "Student "
$tr.Characters(20, 6).Text = "(User)"

# --- Straight Connector 7 (shape id 8): move/resize the lifeline dashed connector
# (target EMU: x=3389242, cx=50713, cy=6807812)
$conn = $s.Shapes.Item(4)
$conn.Left = 266.8695
$conn.Width = 3.9932
$conn.Height = 536.0482
